$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

# "About" sheet: version banner (A2) and recommended-citation text (A6)
# carry the build timestamp.
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A2").Value2 = $wsAbout.Range("A2").Value2.Replace($oldStamp, $newStamp)
$wsAbout.Range("A6").Value2 = $wsAbout.Range("A6").Value2.Replace($oldStamp, $newStamp)

# "Boundaries and methane sources" sheet: build_version column (S),
# one value per data row (rows 2-10).
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")
for ($row = 2; $row -le 10; $row++) {
    $cell = $wsData.Cells.Item($row, 19)  # column S = 19
    $cell.Value2 = $cell.Value2.Replace($oldStamp, $newStamp)
}
